# Ridership run on 20161026.
# Update the Riders (C) and Average (D) columns on the Ridership sheet
# with the refreshed weekly figures. The embedded line chart reads
# directly from these cells, so its cached values will be refreshed
# automatically when Excel recalculates/saves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 229
$ws.Range("D2").Value = 100.26

# Row 3
$ws.Range("C3").Value = 251
$ws.Range("D3").Value = 102.6

# Row 4
$ws.Range("C4").Value = 173
$ws.Range("D4").Value = 107.13

# Row 5
$ws.Range("C5").Value = 193
$ws.Range("D5").Value = 105.94

# Row 6
$ws.Range("C6").Value = 149
$ws.Range("D6").Value = 102.76

# Row 7
$ws.Range("C7").Value = 97
$ws.Range("D7").Value = 48.88

# Row 8
$ws.Range("C8").Value = 77
$ws.Range("D8").Value = 37.45
